# Auto-generated edit script: updates Leve profit-tracking values
# (currentAveragePrice / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to match a refreshed market-board pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 244.77777
$ws.Range("I33").Value = 244.77777
$ws.Range("K33").Value = 244.77777
$ws.Range("M33").Value = -15.77777
$ws.Range("H38").Value = 1072.5
$ws.Range("I38").Value = 145.33333
$ws.Range("J38").Value = 1999.6666
$ws.Range("K38").Value = 435.99999
$ws.Range("L38").Value = 5998.9998
$ws.Range("M38").Value = -63.99998999999997
$ws.Range("N38").Value = -6742.9998
$ws.Range("H43").Value = 3625
$ws.Range("I43").Value = 3625
$ws.Range("K43").Value = 3625
$ws.Range("M43").Value = -3556
$ws.Range("H62").Value = 1500
$ws.Range("I62").Value = 1500
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 1500
$ws.Range("L62").ClearContents()
$ws.Range("N62").ClearContents()
$ws.Range("M62").Value = -876
$ws.Range("H65").Value = 1500
$ws.Range("I65").Value = 1500
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 7500
$ws.Range("L65").ClearContents()
$ws.Range("N65").ClearContents()
$ws.Range("M65").Value = -4380
$ws.Range("H100").Value = 4848.5
$ws.Range("J100").Value = 4799.5
$ws.Range("L100").Value = 4799.5
$ws.Range("N100").Value = -5881.5
$ws.Range("H103").Value = 654.36365
$ws.Range("I103").Value = 339.6
$ws.Range("K103").Value = 1018.8
$ws.Range("M103").Value = -432.8000000000001
$ws.Range("H106").Value = 3500
$ws.Range("I106").Value = 3500
$ws.Range("K106").Value = 3500
$ws.Range("M106").Value = -2869
$ws.Range("H125").Value = 2400
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H127").Value = 926
$ws.Range("I127").Value = 880.3333
$ws.Range("K127").Value = 2640.9999
$ws.Range("M127").Value = 2319.0001
$ws.Range("H129").Value = 2233.3333
$ws.Range("I129").Value = 1850
$ws.Range("K129").Value = 5550
$ws.Range("M129").Value = -550
$ws.Range("H131").Value = 6131
$ws.Range("J131").Value = 5196.5
$ws.Range("L131").Value = 15589.5
$ws.Range("N131").Value = -25669.5
$ws.Range("H138").Value = 3833.8696
$ws.Range("J138").Value = 4008.5715
$ws.Range("L138").Value = 12025.7145
$ws.Range("N138").Value = -22305.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H32").Value = 4378.6816
$ws.Range("I32").Value = 2775.0557
$ws.Range("K32").Value = 2775.0557
$ws.Range("M32").Value = -2488.0557
$ws.Range("H45").Value = 2423.4666
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("H61").Value = 3180.4443
$ws.Range("I61").Value = 2990.5
$ws.Range("K61").Value = 2990.5
$ws.Range("M61").Value = -2778.5
$ws.Range("H74").Value = 3048.1538
$ws.Range("I74").Value = 2885.5
$ws.Range("K74").Value = 2885.5
$ws.Range("M74").Value = -2011.5
$ws.Range("H77").Value = 3048.1538
$ws.Range("I77").Value = 2885.5
$ws.Range("K77").Value = 14427.5
$ws.Range("M77").Value = -10059.5
$ws.Range("H97").Value = 3149.5715
$ws.Range("I97").Value = 2841.1667
$ws.Range("K97").Value = 2841.1667
$ws.Range("M97").Value = -2345.1667
$ws.Range("H102").Value = 2110.9
$ws.Range("I102").Value = 2827.5
$ws.Range("J102").Value = 1633.1666
$ws.Range("K102").Value = 2827.5
$ws.Range("L102").Value = 1633.1666
$ws.Range("M102").Value = -1205.5
$ws.Range("N102").Value = -4877.1666
$ws.Range("H110").Value = 2240.4285
$ws.Range("J110").Value = 3272.5
$ws.Range("L110").Value = 3272.5
$ws.Range("N110").Value = -7362.5
$ws.Range("H122").Value = 15999.75
$ws.Range("I122").Value = 14666.333
$ws.Range("K122").Value = 43998.999
$ws.Range("M122").Value = -41548.999
$ws.Range("H132").Value = 3173.875
$ws.Range("I132").Value = 2913
$ws.Range("K132").Value = 8739
$ws.Range("M132").Value = -6209
$ws.Range("H136").Value = 3180.4443
$ws.Range("I136").Value = 2990.5
$ws.Range("K136").Value = 8971.5
$ws.Range("M136").Value = -6421.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3666.6667
$ws.Range("H36").Value = 199
$ws.Range("I36").Value = 199
$ws.Range("K36").Value = 199
$ws.Range("M36").Value = 335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1333
$ws.Range("J16").Value = 999.5
$ws.Range("L16").Value = 999.5
$ws.Range("N16").Value = -1573.5
$ws.Range("H39").Value = 5525.5
$ws.Range("I39").Value = 1051
$ws.Range("J39").Value = 10000
$ws.Range("K39").Value = 1051
$ws.Range("L39").Value = 10000
$ws.Range("M39").Value = -660
$ws.Range("N39").Value = -10782
$ws.Range("H49").Value = 5525.5
$ws.Range("I49").Value = 1051
$ws.Range("J49").Value = 10000
$ws.Range("K49").Value = 1051
$ws.Range("L49").Value = 10000
$ws.Range("M49").Value = -869
$ws.Range("N49").Value = -10364
$ws.Range("H54").Value = 15000
$ws.Range("J54").Value = 15000
$ws.Range("L54").Value = 15000
$ws.Range("N54").Value = -16316
$ws.Range("H55").Value = 27500
$ws.Range("I55").Value = 5000
$ws.Range("J55").Value = 50000
$ws.Range("K55").Value = 5000
$ws.Range("L55").Value = 50000
$ws.Range("M55").Value = -4685
$ws.Range("N55").Value = -50630
$ws.Range("H86").Value = 7833
$ws.Range("I86").Value = 3749.5
$ws.Range("K86").Value = 3749.5
$ws.Range("M86").Value = -2626.5
$ws.Range("H89").Value = 7833
$ws.Range("I89").Value = 3749.5
$ws.Range("K89").Value = 18747.5
$ws.Range("M89").Value = -13131.5
$ws.Range("H105").Value = 920
$ws.Range("I105").Value = 920
$ws.Range("K105").Value = 920
$ws.Range("M105").Value = 827
$ws.Range("H113").Value = 1333
$ws.Range("J113").Value = 999.5
$ws.Range("L113").Value = 999.5
$ws.Range("N113").Value = -5339.5
$ws.Range("H122").Value = 1125.2
$ws.Range("I122").Value = 1162.7142
$ws.Range("K122").Value = 3488.1426
$ws.Range("M122").Value = -1038.1426
$ws.Range("H132").Value = 5942.2
$ws.Range("I132").Value = 5802.75
$ws.Range("K132").Value = 17408.25
$ws.Range("M132").Value = -14878.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 125133.125
$ws.Range("I4").Value = 152.28572
$ws.Range("K4").Value = 456.85716
$ws.Range("M4").Value = -344.85716
$ws.Range("H140").Value = 3000
$ws.Range("I140").Value = 3000
$ws.Range("K140").Value = 9000
$ws.Range("M140").Value = -3820

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 875
$ws.Range("I113").Value = 750
$ws.Range("K113").Value = 750
$ws.Range("M113").Value = 1420
$ws.Range("H122").Value = 13061.75
$ws.Range("I122").Value = 13061.75
$ws.Range("K122").Value = 39185.25
$ws.Range("M122").Value = -36735.25
$ws.Range("H132").Value = 3649.5334
$ws.Range("I132").Value = 3700.3845
$ws.Range("J132").Value = 3319
$ws.Range("K132").Value = 11101.1535
$ws.Range("L132").Value = 9957
$ws.Range("M132").Value = -8571.1535
$ws.Range("N132").Value = -15017

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6008.5
$ws.Range("I7").Value = 5020.6665
$ws.Range("J7").Value = 6749.375
$ws.Range("K7").Value = 5020.6665
$ws.Range("L7").Value = 6749.375
$ws.Range("M7").Value = -4908.6665
$ws.Range("N7").Value = -6973.375
$ws.Range("H22").Value = 666.1667
$ws.Range("I22").Value = 499.25
$ws.Range("K22").Value = 499.25
$ws.Range("M22").Value = -204.25
$ws.Range("H27").Value = 666.1667
$ws.Range("I27").Value = 499.25
$ws.Range("K27").Value = 499.25
$ws.Range("M27").Value = -392.25
$ws.Range("H126").Value = 6008.5
$ws.Range("I126").Value = 5020.6665
$ws.Range("J126").Value = 6749.375
$ws.Range("K126").Value = 15061.9995
$ws.Range("L126").Value = 20248.125
$ws.Range("M126").Value = -12591.9995
$ws.Range("N126").Value = -25188.125
$ws.Range("H132").Value = 16283.077
$ws.Range("I132").Value = 18408.363
$ws.Range("J132").Value = 4594
$ws.Range("K132").Value = 55225.08900000001
$ws.Range("L132").Value = 13782
$ws.Range("M132").Value = -52695.08900000001
$ws.Range("N132").Value = -18842

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 30000
$ws.Range("I40").Value = 30000
$ws.Range("K40").Value = 30000
$ws.Range("M40").Value = -29851
$ws.Range("H100").Value = 20999.5
$ws.Range("I100").Value = 20999.5
$ws.Range("K100").Value = 41999
$ws.Range("M100").Value = -41458
$ws.Range("H132").Value = 4000.25
$ws.Range("I132").Value = 4000.25
$ws.Range("K132").Value = 12000.75
$ws.Range("M132").Value = -9470.75

